$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price observation was recorded for this market/product.
# Insert a new row at row 65 (pushing the existing row 65..155 data down
# to rows 66..156) and populate it with the new observation.
$ws.Rows(65).Insert()

$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 45079
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112022
$ws.Range("G65").Value = "Arveja Verde"
$ws.Range("H65").Value = "Perfection"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 40
$ws.Range("K65").Value = 43000
$ws.Range("L65").Value = 43000
$ws.Range("M65").Value = 43000
$ws.Range("N65").Value = "$/malla 25 kilos"
$ws.Range("O65").Value = "Provincia de Huasco"
$ws.Range("P65").Value = 1720
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = "Hortaliza"
